$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove row 260 (colo "XNN" - Xining, China). This shifts all subsequent
# rows up by one, matching the target dimension A1:H330.
$ws.Rows.Item(260).Delete()
